$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update / insert the Cluster Name values (column A), rows 2-12
$ws.Range("A2").Value = "3516 The Alexander Aged Care Centre"
$ws.Range("A3").Value = "3641 Calvary The Regent Mount Waverley"
$ws.Range("A4").Value = "Barwon Heads Hotel Barwon Heads"
$ws.Range("A5").Value = "Confirmed Omicron Sircuit Bar Fitzroy"
$ws.Range("A6").Value = "Confirmed Omicron Variant The Peel Hotel Collingwood"
$ws.Range("A7").Value = "Diamond Valley Pork and Baxters Pork Laverton North"
$ws.Range("A8").Value = "Melbourne Cricket Ground (MCG)"
$ws.Range("A9").Value = "Melbourne Stars Big Bash Cricket Team EastMelbourne"
$ws.Range("A10").Value = "The Royal Children's Hospital Melbourne Emergency Department Parkville"
$ws.Range("A11").Value = "Werribee Mercy Hospital Emergency Department"
$ws.Range("A12").Value = "Western Health Sunshine Hospital Emergency Department St Albans"

# Update the Activecases values (column B), rows 2-12
$ws.Range("B2").Value = 14
$ws.Range("B3").Value = 10
$ws.Range("B4").Value = 16
$ws.Range("B5").Value = 18
$ws.Range("B6").Value = 14
$ws.Range("B7").Value = 10
$ws.Range("B8").Value = 47
$ws.Range("B9").Value = 24
$ws.Range("B10").Value = 10
$ws.Range("B11").Value = 13
$ws.Range("B12").Value = 11
